$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "305.06"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.10%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.86"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.92%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.979"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.65%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08082"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.05%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.901"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.16%"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.73%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.877"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.58%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9301"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.34%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1244"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-18.03%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1905"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.31%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09215"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.21%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03528"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.64%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09928"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.53%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001413"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.17%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006052"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "4.43%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.599"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.68%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.105"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "3.85%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3451"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.22%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.223"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "3.94%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1294"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.90%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2532"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "6.05%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04414"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.28%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001234"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2.51%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004713"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-3.18%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001301"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "6.33%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003131"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-29.14%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01953"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.82%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05254"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "8.98%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007541"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.74%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01013"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.53%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1373"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.37%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002102"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.91%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01072"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.42%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006346"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "4.22%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.61%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.57"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-1.70%"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "40.19%"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.61%"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.61%"
